$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bulk price/volume updates for rows 2-49 (coin identity unchanged)
# D-column values are prefixed with a leading apostrophe so Excel keeps
# them as text (matching the original "inline string" price formatting)
# instead of auto-converting to a numeric type.
$ws.Range("D2").Value = "'72.823.51"
$ws.Range("E2").Value = "  +1.18%  "
$ws.Range("D3").Value = "'3.982.40"
$ws.Range("E3").Value = "  -0.38%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'585.68"
$ws.Range("E5").Value = "  +8.11%  "
$ws.Range("D6").Value = "'158.40"
$ws.Range("E6").Value = "  +5.12%  "
$ws.Range("D7").Value = "'0.680"
$ws.Range("E7").Value = "  -2.94%  "
$ws.Range("E8").Value = "  -0.19%  "
$ws.Range("D9").Value = "'0.751"
$ws.Range("E9").Value = "  +0.47%  "
$ws.Range("D10").Value = "'0.168"
$ws.Range("E10").Value = "  -1.54%  "
$ws.Range("D11").Value = "'54.28"
$ws.Range("E11").Value = "  +1.45%  "
$ws.Range("D12").Value = "'0.0000318"
$ws.Range("E12").Value = "  -1.88%  "
$ws.Range("D13").Value = "'10.87"
$ws.Range("E13").Value = "  +1.78%  "
$ws.Range("D14").Value = "'4.612.01"
$ws.Range("E14").Value = "  -0.55%  "
$ws.Range("D15").Value = "'3.967.13"
$ws.Range("E15").Value = "  -0.55%  "
$ws.Range("D16").Value = "'1.28"
$ws.Range("E16").Value = "  +7.96%  "
$ws.Range("D17").Value = "'14.03"
$ws.Range("E17").Value = "  -0.76%  "
$ws.Range("D18").Value = "'20.47"
$ws.Range("E18").Value = "  -0.44%  "
$ws.Range("D20").Value = "'72.447.51"
$ws.Range("E20").Value = "  +0.82%  "
$ws.Range("D21").Value = "'433.91"
$ws.Range("E21").Value = "  +0.17%  "
$ws.Range("D22").Value = "'4.69"
$ws.Range("E22").Value = "  +9.80%  "
$ws.Range("D23").Value = "'96.04"
$ws.Range("E23").Value = "  -0.81%  "
$ws.Range("D24").Value = "'3.43"
$ws.Range("E24").Value = "  -3.45%  "
$ws.Range("D25").Value = "'14.30"
$ws.Range("E25").Value = "  -0.48%  "
$ws.Range("D26").Value = "'4.40"
$ws.Range("E26").Value = "  +21.42%  "
$ws.Range("E27").Value = "  -2.31%  "
$ws.Range("E28").Value = "  +0.38%  "
$ws.Range("D29").Value = "'5.93"
$ws.Range("E29").Value = "  +1.13%  "
$ws.Range("D30").Value = "'36.42"
$ws.Range("E30").Value = "  -0.94%  "
$ws.Range("D31").Value = "'7.83"
$ws.Range("E31").Value = "  +4.13%  "
$ws.Range("D32").Value = "'50.76"
$ws.Range("E32").Value = "  +3.47%  "
$ws.Range("D33").Value = "'13.63"
$ws.Range("E33").Value = "  +1.42%  "
$ws.Range("E34").Value = "  -0.06%  "
$ws.Range("D35").Value = "'678.30"
$ws.Range("E35").Value = "  -0.18%  "
$ws.Range("D36").Value = "'68.74"
$ws.Range("E36").Value = "  +3.97%  "
$ws.Range("E37").Value = "  -2.28%  "
$ws.Range("D38").Value = "'0.0₃0859"
$ws.Range("E38").Value = "  +3.43%  "
$ws.Range("D39").Value = "'3.38"
$ws.Range("E39").Value = "  +0.82%  "
$ws.Range("D40").Value = "'0.998"
$ws.Range("E40").Value = "  -0.19%  "
$ws.Range("D41").Value = "'0.146"
$ws.Range("E41").Value = "  -4.36%  "
$ws.Range("D42").Value = "'3.33"
$ws.Range("E42").Value = "  -2.22%  "
$ws.Range("E43").Value = "  +0.01%  "
$ws.Range("E44").Value = "  +11.30%  "
$ws.Range("D45").Value = "'0.0487"
$ws.Range("E45").Value = "  -0.21%  "
$ws.Range("D46").Value = "'0.149"
$ws.Range("E46").Value = "  -0.78%  "
$ws.Range("D47").Value = "'2.69"
$ws.Range("E47").Value = "  -2.74%  "
$ws.Range("D48").Value = "'3.37"
$ws.Range("E48").Value = "  -0.08%  "
$ws.Range("D49").Value = "'3.43"

# Rows 50-51: Stacks and ARBITRUM swap rank order, with new price/volume values
$ws.Range("B50").Value = "ARBITRUM"
$ws.Range("C50").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D50").Value = "'2.16"
$ws.Range("E50").Value = "  +7.49%  "

$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").Value = "'3.00"
$ws.Range("E51").Value = "  +0.17%  "
